# Updated cryptos list on Wed Feb 21 18:12:13 UTC 2024 with GitHub Actions
#
# Refreshes the "cryptos" worksheet with a new coinranking.com snapshot:
#   - Price (col D) and Volume(1h) (col E) are refreshed for every coin
#     that is still present in the table.
#   - Rows 26-51 (ranks 24-49) shift up by one place: each coin's new
#     row now shows the Name/Link/Price/Volume that the coin one rank
#     below it used to show.
#   - A brand-new coin (BEAM) enters the table at the bottom, row 51.
#
# All of the table's Price/Volume cells are stored as plain text in the
# workbook (e.g. "50.941.20", "  -0.80%  "), so each write forces a text
# number-format before assigning the value and then clears the (blank)
# formatting again, to avoid Excel "helpfully" re-interpreting numeric-
# looking strings (like "366.76" or "1.00") as actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '50.941.20' },
    @{ Cell = 'E2'; Value = '  -0.80%  ' },
    @{ Cell = 'D3'; Value = '2.899.72' },
    @{ Cell = 'E3'; Value = '  -0.54%  ' },
    @{ Cell = 'E4'; Value = '  -0.19%  ' },
    @{ Cell = 'D5'; Value = '366.76' },
    @{ Cell = 'E5'; Value = '  +5.00%  ' },
    @{ Cell = 'D6'; Value = '102.41' },
    @{ Cell = 'E6'; Value = '  -3.00%  ' },
    @{ Cell = 'D7'; Value = '0.538' },
    @{ Cell = 'E7'; Value = '  -2.77%  ' },
    @{ Cell = 'D8'; Value = '1.00' },
    @{ Cell = 'E8'; Value = '  -0.06%  ' },
    @{ Cell = 'D9'; Value = '0.582' },
    @{ Cell = 'E9'; Value = '  -3.61%  ' },
    @{ Cell = 'D10'; Value = '36.32' },
    @{ Cell = 'E10'; Value = '  -3.40%  ' },
    @{ Cell = 'E11'; Value = '  +0.58%  ' },
    @{ Cell = 'D12'; Value = '0.0828' },
    @{ Cell = 'E12'; Value = '  -1.99%  ' },
    @{ Cell = 'E13'; Value = '  -3.50%  ' },
    @{ Cell = 'D14'; Value = '3.349.73' },
    @{ Cell = 'E14'; Value = '  -0.64%  ' },
    @{ Cell = 'D15'; Value = '7.33' },
    @{ Cell = 'E15'; Value = '  -3.26%  ' },
    @{ Cell = 'D16'; Value = '2.893.98' },
    @{ Cell = 'E16'; Value = '  -0.58%  ' },
    @{ Cell = 'D17'; Value = '0.921' },
    @{ Cell = 'E17'; Value = '  -3.81%  ' },
    @{ Cell = 'D18'; Value = '50.865.51' },
    @{ Cell = 'E18'; Value = '  -0.92%  ' },
    @{ Cell = 'D19'; Value = '3.20' },
    @{ Cell = 'E19'; Value = '  -7.05%  ' },
    @{ Cell = 'D20'; Value = '7.12' },
    @{ Cell = 'E20'; Value = '  -3.56%  ' },
    @{ Cell = 'D21'; Value = '12.79' },
    @{ Cell = 'E21'; Value = '  -4.20%  ' },
    @{ Cell = 'D22'; Value = '0.0₃0938' },
    @{ Cell = 'E22'; Value = '  -2.55%  ' },
    @{ Cell = 'D23'; Value = '67.92' },
    @{ Cell = 'E23'; Value = '  -1.17%  ' },
    @{ Cell = 'D24'; Value = '257.72' },
    @{ Cell = 'E24'; Value = '  -0.51%  ' },
    @{ Cell = 'E25'; Value = '  -0.72%  ' },
    @{ Cell = 'B26'; Value = 'Kaspa' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' },
    @{ Cell = 'D26'; Value = '0.174' },
    @{ Cell = 'E26'; Value = '  +0.39%  ' },
    @{ Cell = 'B27'; Value = 'Dai' },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D27'; Value = '1.00' },
    @{ Cell = 'E27'; Value = '  +0.06%  ' },
    @{ Cell = 'B28'; Value = 'EthereumClassic' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Cell = 'D28'; Value = '25.47' },
    @{ Cell = 'E28'; Value = '  -3.10%  ' },
    @{ Cell = 'B29'; Value = 'Filecoin' },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D29'; Value = '6.93' },
    @{ Cell = 'E29'; Value = '  -6.21%  ' },
    @{ Cell = 'B30'; Value = 'Hedera' },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Cell = 'D30'; Value = '0.101' },
    @{ Cell = 'E30'; Value = '  -2.41%  ' },
    @{ Cell = 'B31'; Value = 'RenderToken' },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Cell = 'D31'; Value = '6.16' },
    @{ Cell = 'E31'; Value = '  +1.37%  ' },
    @{ Cell = 'B32'; Value = 'Cosmos' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Cell = 'D32'; Value = '9.84' },
    @{ Cell = 'E32'; Value = '  -3.46%  ' },
    @{ Cell = 'B33'; Value = 'Toncoin' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Cell = 'D33'; Value = '2.12' },
    @{ Cell = 'E33'; Value = '  -1.30%  ' },
    @{ Cell = 'B34'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D34'; Value = '34.19' },
    @{ Cell = 'E34'; Value = '  -3.64%  ' },
    @{ Cell = 'B35'; Value = 'OKB' },
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' },
    @{ Cell = 'D35'; Value = '50.83' },
    @{ Cell = 'E35'; Value = '  +0.95%  ' },
    @{ Cell = 'B36'; Value = 'FirstDigitalUSD' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd' },
    @{ Cell = 'D36'; Value = '1.00' },
    @{ Cell = 'E36'; Value = '  +0.42%  ' },
    @{ Cell = 'B37'; Value = 'VeChain' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D37'; Value = '0.0419' },
    @{ Cell = 'E37'; Value = '  -1.09%  ' },
    @{ Cell = 'B38'; Value = 'LidoDAOToken' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Cell = 'D38'; Value = '2.97' },
    @{ Cell = 'E38'; Value = '  -4.65%  ' },
    @{ Cell = 'B39'; Value = 'Stacks' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' },
    @{ Cell = 'D39'; Value = '2.61' },
    @{ Cell = 'E39'; Value = '  -0.86%  ' },
    @{ Cell = 'B40'; Value = 'Celestia' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia' },
    @{ Cell = 'D40'; Value = '16.89' },
    @{ Cell = 'E40'; Value = '  -3.88%  ' },
    @{ Cell = 'B41'; Value = 'ARBITRUM' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Cell = 'D41'; Value = '1.82' },
    @{ Cell = 'E41'; Value = '  -5.21%  ' },
    @{ Cell = 'B42'; Value = 'Stellar' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Cell = 'D42'; Value = '0.112' },
    @{ Cell = 'E42'; Value = '  -2.93%  ' },
    @{ Cell = 'B43'; Value = 'EnergySwap' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D43'; Value = '21.90' },
    @{ Cell = 'E43'; Value = '  -1.61%  ' },
    @{ Cell = 'B44'; Value = 'Monero' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D44'; Value = '118.08' },
    @{ Cell = 'E44'; Value = '  -1.61%  ' },
    @{ Cell = 'B45'; Value = 'WEMIXToken' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D45'; Value = '2.07' },
    @{ Cell = 'E45'; Value = '  -2.40%  ' },
    @{ Cell = 'B46'; Value = 'Maker' },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D46'; Value = '2.007.92' },
    @{ Cell = 'E46'; Value = '  -4.05%  ' },
    @{ Cell = 'B47'; Value = 'ApeXProtocol' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' },
    @{ Cell = 'D47'; Value = '2.31' },
    @{ Cell = 'E47'; Value = '  -0.38%  ' },
    @{ Cell = 'B48'; Value = 'NEARProtocol' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Cell = 'D48'; Value = '3.14' },
    @{ Cell = 'E48'; Value = '  -4.54%  ' },
    @{ Cell = 'B49'; Value = 'RocketPoolETH' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' },
    @{ Cell = 'D49'; Value = '3.190.32' },
    @{ Cell = 'E49'; Value = '  -0.30%  ' },
    @{ Cell = 'B50'; Value = 'TheGraph' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' },
    @{ Cell = 'D50'; Value = '0.236' },
    @{ Cell = 'E50'; Value = '  -0.62%  ' },
    @{ Cell = 'B51'; Value = 'BEAM' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam' },
    @{ Cell = 'D51'; Value = '0.0307' },
    @{ Cell = 'E51'; Value = '  -7.51%  ' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.ClearFormats()
}
